# Adds two new IPT coverage data rows (row 6: unit cost, row 7: total cost)
# to the "time_variants" worksheet, mirroring the existing vaccination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: econ_program_unitcost_ipt ---
$ws.Range("A6").Value = "econ_program_unitcost_ipt"
$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "yes"

# Columns E6:AX6 are all zero.
$zeroCols6 = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
               "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR", `
               "AS","AT","AU","AV","AW","AX")
foreach ($col in $zeroCols6) {
    $ws.Range($col + "6").Value = 0
}

$ws.Range("AY6").Value = 20.77
$ws.Range("BB6").Value = 22
$ws.Range("BE6").Value = 21
$ws.Range("BF6").Value = 21
$ws.Range("BG6").Value = 21
$ws.Range("BH6").Value = 21
$ws.Range("BI6").Value = 21

# --- Row 7: econ_program_totalcost_ipt ---
$ws.Range("A7").Value = "econ_program_totalcost_ipt"
$ws.Range("B7").Value = "yes"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "yes"

$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("Q7").Value = 30000
$ws.Range("V7").Value = 50000
$ws.Range("AF7").Value = 80000
$ws.Range("AK7").Value = 100000
$ws.Range("AP7").Value = 100000
$ws.Range("AZ7").Value = 100000
$ws.Range("BE7").Value = 110000
$ws.Range("BF7").Value = 250000
$ws.Range("BH7").Value = 250000
$ws.Range("BI7").Value = 250000

# Update the active selection/view to match the authored workbook state.
$ws.Range("A8").Select()
